# edit.ps1 - applies the ЛР4_ОП.docx diff via Word COM-interop
$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument
$wdFindContinue = 1

function Find-ParagraphByText {
    param([string]$searchText)
    $rng = $d.Content
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
    if (-not $ok) {
        throw ("Could not find paragraph containing: " + $searchText)
    }
    return $rng.Paragraphs(1)
}

# --- 1) Remove the empty, italic-styled paragraph that sits right before
#        the "Виконання мовою Python." heading (two paragraphs back: a
#        "tabs only" paragraph, then the empty italic one). ---
$pPython = Find-ParagraphByText("Python.")
$pTabs = $pPython.Previous()
$pEmptyItalic = $pTabs.Previous()
[void]$pEmptyItalic.Range.Delete()

# --- 2) n=int(input("...")): wrap "int" in spellStart/spellEnd proofErr
#        marks, and split the "(input(" run into "(" + input (wrapped) + "(" ---
$xmlNInt = @'
<w:p w14:paraId="5386173D" w14:textId="77777777" w:rsidR="006B1ADF" w:rsidRPr="006B1ADF" w:rsidRDefault="006B1ADF" w:rsidP="006B1ADF" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>n=</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="2B91AF"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>input</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>(</w:t></w:r><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>"Введіть кількість членів n: "</w:t></w:r><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">)) </w:t></w:r></w:p>
'@
$pNInt = Find-ParagraphByText("n=int(input(")
[void]$pNInt.Range.InsertXML($xmlNInt)

# --- 3) for i in range(n-1): wrap for/in/range with proofErr marks ---
$xmlFor = @'
<w:p w14:paraId="7FD673AC" w14:textId="77777777" w:rsidR="006B1ADF" w:rsidRPr="006B1ADF" w:rsidRDefault="006B1ADF" w:rsidP="006B1ADF" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>for</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> i </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>in</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="2B91AF"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>range</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">(n-1): </w:t></w:r><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="008000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>#Оскільки перший член уже існує, то потрібно порахувати ще (n-1) членів</w:t></w:r></w:p>
'@
$pFor = Find-ParagraphByText("for i in range")
[void]$pFor.Range.InsertXML($xmlFor)

# --- 4) Comment run split around "переприсвоюємо" with proofErr marks ---
$xmlComment = @'
<w:p w14:paraId="612077CE" w14:textId="77777777" w:rsidR="006B1ADF" w:rsidRPr="006B1ADF" w:rsidRDefault="006B1ADF" w:rsidP="006B1ADF" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">    a=(a**2)/(a+3) </w:t></w:r><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="008000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">#Шукаємо наступний член за формулою і </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="008000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>переприсвоюємо</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="008000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> значення змінної члена</w:t></w:r></w:p>
'@
$pComment = Find-ParagraphByText("a=(a**2)/(a+3)")
[void]$pComment.Range.InsertXML($xmlComment)

# --- 5) print(...) split into "print" (wrapped) + "(" ---
$xmlPrint = @'
<w:p w14:paraId="7F00D3DB" w14:textId="77777777" w:rsidR="006B1ADF" w:rsidRPr="006B1ADF" w:rsidRDefault="006B1ADF" w:rsidP="006B1ADF" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>print</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>(</w:t></w:r><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>"Сума елементів послідовности ="</w:t></w:r><w:r w:rsidRPr="006B1ADF"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>,s)</w:t></w:r></w:p>
'@
$pPrint = Find-ParagraphByText(",s)")
[void]$pPrint.Range.InsertXML($xmlPrint)

# --- 6) Insert a new empty paragraph right after the print(...) paragraph
#        (and before the existing empty paragraph that was already there). ---
$xmlInsertCombo = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p w14:paraId="62045478" w14:textId="77777777" w:rsidR="006B1ADF" w:rsidRPr="00774095" w:rsidRDefault="006B1ADF" w:rsidP="005B4495" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p>
'@
$pPrint2 = Find-ParagraphByText(",s)")
$pEmptyAfterPrint = $pPrint2.Next()
[void]$pEmptyAfterPrint.Range.InsertXML($xmlInsertCombo)

Write-Output "edit.ps1 completed successfully"
